{"js": "// Update the header date and the twenty-five two-digit multiplication\n// prompts in the table, replacing each old value with its new value.\n// Each old string is unique in the document, so an exact, case-sensitive\n// search-and-replace per pair is safe and keeps existing run formatting.\nconst replacements = [\n  [\"2025-08-14 Thursday\", \"2025-08-15 Friday\"],\n  [\"42\u00d779=\", \"81\u00d725=\"],\n  [\"85\u00d722=\", \"60\u00d727=\"],\n  [\"28\u00d761=\", \"94\u00d721=\"],\n  [\"65\u00d754=\", \"71\u00d759=\"],\n  [\"78\u00d721=\", \"87\u00d723=\"],\n  [\"67\u00d792=\", \"11\u00d713=\"],\n  [\"66\u00d746=\", \"63\u00d711=\"],\n  [\"68\u00d745=\", \"92\u00d779=\"],\n  [\"78\u00d757=\", \"94\u00d762=\"],\n  [\"93\u00d724=\", \"87\u00d732=\"],\n  [\"72\u00d749=\", \"66\u00d747=\"],\n  [\"46\u00d774=\", \"94\u00d734=\"],\n  [\"29\u00d743=\", \"36\u00d798=\"],\n  [\"28\u00d750=\", \"36\u00d795=\"],\n  [\"72\u00d796=\", \"23\u00d739=\"],\n  [\"27\u00d729=\", \"63\u00d716=\"],\n  [\"27\u00d718=\", \"77\u00d777=\"],\n  [\"63\u00d745=\", \"67\u00d791=\"],\n  [\"89\u00d777=\", \"60\u00d765=\"],\n  [\"69\u00d782=\", \"46\u00d757=\"],\n  [\"33\u00d760=\", \"90\u00d796=\"],\n  [\"59\u00d763=\", \"63\u00d721=\"],\n  [\"72\u00d742=\", \"71\u00d774=\"],\n  [\"12\u00d777=\", \"21\u00d761=\"],\n  [\"64\u00d732=\", \"72\u00d721=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the header date and the twenty-five two-digit multiplication\n# prompts in the table, replacing each old value with its new value.\n# Each old string is unique in the document, so a simple Find/Replace\n# per pair (case-sensitive, whole document) is safe and preserves the\n# existing run formatting.\n$d = $word.ActiveDocument\n\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\n$replacements = @(\n  @(\"2025-08-14 Thursday\", \"2025-08-15 Friday\"),\n  @(\"42\u00d779=\", \"81\u00d725=\"),\n  @(\"85\u00d722=\", \"60\u00d727=\"),\n  @(\"28\u00d761=\", \"94\u00d721=\"),\n  @(\"65\u00d754=\", \"71\u00d759=\"),\n  @(\"78\u00d721=\", \"87\u00d723=\"),\n  @(\"67\u00d792=\", \"11\u00d713=\"),\n  @(\"66\u00d746=\", \"63\u00d711=\"),\n  @(\"68\u00d745=\", \"92\u00d779=\"),\n  @(\"78\u00d757=\", \"94\u00d762=\"),\n  @(\"93\u00d724=\", \"87\u00d732=\"),\n  @(\"72\u00d749=\", \"66\u00d747=\"),\n  @(\"46\u00d774=\", \"94\u00d734=\"),\n  @(\"29\u00d743=\", \"36\u00d798=\"),\n  @(\"28\u00d750=\", \"36\u00d795=\"),\n  @(\"72\u00d796=\", \"23\u00d739=\"),\n  @(\"27\u00d729=\", \"63\u00d716=\"),\n  @(\"27\u00d718=\", \"77\u00d777=\"),\n  @(\"63\u00d745=\", \"67\u00d791=\"),\n  @(\"89\u00d777=\", \"60\u00d765=\"),\n  @(\"69\u00d782=\", \"46\u00d757=\"),\n  @(\"33\u00d760=\", \"90\u00d796=\"),\n  @(\"59\u00d763=\", \"63\u00d721=\"),\n  @(\"72\u00d742=\", \"71\u00d774=\"),\n  @(\"12\u00d777=\", \"21\u00d761=\"),\n  @(\"64\u00d732=\", \"72\u00d721=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($oldText, $true, $true, $false, $false, $false, $true, $wdFindContinue, $false, $newText, $wdReplaceAll) | Out-Null\n}\n"}
